$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the answer explanation placeholder for every question row (26-100)
# by filling column G with a placeholder string, matching the style already
# used by the surrounding cells (font color FF333333, i.e. existing style 1).
$rng = $ws.Range("G26:G100")
$rng.Value = "I'm sure you know why (Placeholder)"
$rng.Font.Color = 3355443

# Remember the selection/scroll position at time of submit: the active cell
# moves from G100 to G105 (off the used range, since the sheet view just
# tracks the last click position).
[void]$ws.Range("G105").Select()
